$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.924.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0686"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.93%  "
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.781.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.910.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  -2.62%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0551"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  +8.96%  "
$ws.Range("E35").Value = "  +6.91%  "
$ws.Range("E36").Value = "  +14.75%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.348.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.022.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +22.13%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0665"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
